# Daily attendance processing - 2026-01-10 12:48:20
# Normalize the "Recorded By" (column G) lists: move the exact token "System"
# (case-sensitive) to the front of the comma-separated list, keeping the
# remaining entries in their original relative order. If the exact token
# "System" is not present, sort the entries alphabetically instead.

function Test-CaseSensitiveEquals($s1, $s2) {
    if ($s1.Length -ne $s2.Length) { return $false }
    for ($i = 0; $i -lt $s1.Length; $i++) {
        if ([int]$s1[$i] -ne [int]$s2[$i]) { return $false }
    }
    return $true
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val -or $val -eq "") {
        continue
    }

    $parts = $val -split ", "

    if ($parts.Count -le 1) {
        continue
    }

    $hasSystem = $false
    foreach ($p in $parts) {
        if (Test-CaseSensitiveEquals $p "System") {
            $hasSystem = $true
        }
    }

    if ($hasSystem) {
        $rest = @()
        foreach ($p in $parts) {
            if (-not (Test-CaseSensitiveEquals $p "System")) {
                $rest += $p
            }
        }
        $newParts = @("System") + $rest
    } else {
        $newParts = $parts | Sort-Object
    }

    $newVal = $newParts -join ", "

    if ($newVal -ne $val) {
        $cell.Value2 = $newVal
    }
}
